$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha2")

# --- Fill in the ID column (A) for the requirements that already exist ---
$ids = @("RF01","RF02","RF03","RF04","RF05","RF06","RF07","RF08","RF09","RF10","RF11","RF12","RF13","RF14","RNF01","RNF02","RNF03","RNF04","RNF05","RNF06","RNF07","RNF08")
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $ids[$i]
}

# --- Insert the new requirement row (RF15) before the current row 16 ---
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = "RF15"

# --- Fix existing text / priority values (typos & reclassification) ---
$ws.Range("B3").Value  = "O software deve realizar autenticação de login e senha"
$ws.Range("C4").Value  = "Importante"
$ws.Range("C5").Value  = "Importante"
$ws.Range("C6").Value  = "Importante"
$ws.Range("C9").Value  = "Importante"
$ws.Range("B10").Value = "O software deve ter opção de seleção de jogos"
$ws.Range("C12").Value = "Essencial"
$ws.Range("C13").Value = "Desejável"
$ws.Range("C14").Value = "Desejável"
$ws.Range("C15").Value = "Desejável"

$ws.Range("B16").Value = "O software deve mostrar os próximos jogos na tela home"
$ws.Range("C16").Value = "Essencial"
$ws.Range("D16").Value = "Funcional "

# --- Append a brand-new requirement row at the end (row 25) ---
$ws.Range("A25").Value = "RF16"
$ws.Range("B25").Value = "O projeto deve ter um diagrama de classe"
$ws.Range("C25").Value = "Essencial"
$ws.Range("D25").Value = "Funcional "
$ws.Range("E25").Value = 43891
$ws.Range("E25").NumberFormat = "m/d/yy"
$ws.Range("G25").Value = "Gustavo Henrique "

# --- Turn the requirements range into a proper Excel Table ---
$rng = $ws.Range("A1:I25")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Tabela1"
$tbl.TableStyle = "TableStyleMedium2"

# --- Hide the audit-trail columns, keep only ID/Descrição/Prioridade/Tipo visible ---
$ws.Range("E1:I1").EntireColumn.Hidden = $true

$ws.Range("J19").Select()
